{"js": "// Append the new \"TODO\" section to the end of the document body.\n// The existing body content (prospect description) is left untouched;\n// we only add two blank paragraphs followed by the TODO heading and its\n// list items, plus a trailing blank paragraph, matching the target\n// OOXML produced by the commit \"bg and favicon added\".\n\nconst body = context.document.body;\n\nconst newParagraphs = [\n  \"\",\n  \"\",\n  \"TODO\",\n  \"Add favicon\",\n  \"Add bg photo,  delete stock pics\",\n  \"Create model of cafes\",\n  \"Add caf\u00e9 to DB\",\n  \"Create superuser\",\n  \"Create form for adding new places\",\n  \"\",\n];\n\nfor (const text of newParagraphs) {\n  body.insertParagraph(text, Word.InsertLocation.end);\n}\n\nawait context.sync();\n", "ps1": "# Append the new \"TODO\" section to the end of the document body.\n# The existing body content (prospect description) is left untouched;\n# we only add two blank paragraphs followed by the TODO heading and its\n# list items, plus a trailing blank paragraph, matching the target\n# OOXML produced by the commit \"bg and favicon added\".\n\n$d = $word.ActiveDocument\n\n$newParagraphs = @(\n    \"\",\n    \"\",\n    \"TODO\",\n    \"Add favicon\",\n    \"Add bg photo,  delete stock pics\",\n    \"Create model of cafes\",\n    \"Add caf\u00e9 to DB\",\n    \"Create superuser\",\n    \"Create form for adding new places\",\n    \"\"\n)\n\nforeach ($text in $newParagraphs) {\n    $endRange = $d.Content\n    $endRange.Collapse(0)              # wdCollapseEnd\n    $endRange.InsertParagraphAfter()   # adds a new, empty trailing paragraph\n\n    $newPara = $d.Content\n    $newPara.Collapse(0)               # move to the freshly created paragraph\n    if ($text -ne \"\") {\n        $newPara.Text = $text\n    }\n}\n"}
